# Weekly update: insert a new price record at the top of the Coliflor /
# Macroferia Regional de Talca data block (row 458), pushing the existing
# rows 458:504 down to 459:505.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 458 (shifts 458:504 down to 459:505, extends
# dimension from A1:R504 to A1:R505).
$ws.Rows(458).Insert()

# Populate the new row 458 with this week's record.
$ws.Cells.Item(458, 1).Value = 5
$ws.Cells.Item(458, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(458, 3).Value = "Maule"
$ws.Cells.Item(458, 4).Value = 45194
$ws.Cells.Item(458, 5).Value = 7
$ws.Cells.Item(458, 6).Value = 100112008
$ws.Cells.Item(458, 7).Value = "Coliflor"
$ws.Cells.Item(458, 8).Value = "Sin especificar"
$ws.Cells.Item(458, 9).Value = "Primera"
$ws.Cells.Item(458, 10).Value = 3000
$ws.Cells.Item(458, 11).Value = 900
$ws.Cells.Item(458, 12).Value = 900
$ws.Cells.Item(458, 13).Value = 900
$ws.Cells.Item(458, 14).Value = "`$/unidad"
$ws.Cells.Item(458, 15).Value = "Región del Maule"
$ws.Cells.Item(458, 16).Value = 900
$ws.Cells.Item(458, 17).Value = 1
$ws.Cells.Item(458, 18).Value = "Hortaliza"
